$wb = $excel.ActiveWorkbook

# --- Sheet ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H19").Value = 2265
$ws.Range("I19").Value = 1243
$ws.Range("K19").Value = 1243
$ws.Range("M19").Value = -1068

$ws.Range("H53").Value = 249.2
$ws.Range("I53").Value = 106.44444
$ws.Range("J53").Value = 463.33334
$ws.Range("K53").Value = 106.44444
$ws.Range("L53").Value = 463.33334
$ws.Range("M53").Value = 530.55556
$ws.Range("N53").Value = -1737.33334

$ws.Range("H97").Value = 2021.1
$ws.Range("J97").Value = 2212.3333
$ws.Range("L97").Value = 6636.999899999999
$ws.Range("N97").Value = -7628.999899999999

$ws.Range("H98").Value = 504
$ws.Range("I98").Value = 397.2143
$ws.Range("K98").Value = 397.2143
$ws.Range("M98").Value = 1100.7857

$ws.Range("H100").Value = 1713.2354
$ws.Range("J100").Value = 3200.8333
$ws.Range("L100").Value = 3200.8333
$ws.Range("N100").Value = -4282.8333

$ws.Range("H107").Value = 864.619
$ws.Range("I107").Value = 600.5625
$ws.Range("K107").Value = 600.5625
$ws.Range("M107").Value = 1319.4375

$ws.Range("H110").Value = 59991
$ws.Range("J110").Value = 59991
$ws.Range("L110").Value = 59991
$ws.Range("N110").Value = -68171

$ws.Range("H122").Value = 504
$ws.Range("I122").Value = 397.2143
$ws.Range("K122").Value = 1191.6429
$ws.Range("M122").Value = 1258.3571

$ws.Range("H132").Value = 1948.8928
$ws.Range("I132").Value = 1596.875
$ws.Range("K132").Value = 4790.625
$ws.Range("M132").Value = -2260.625

# --- Sheet ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1728
$ws.Range("I2").Value = 1449.5
$ws.Range("K2").Value = 1449.5
$ws.Range("M2").Value = -1336.5

$ws.Range("H32").Value = 21746540
$ws.Range("I32").Value = 27031414
$ws.Range("K32").Value = 27031414
$ws.Range("M32").Value = -27031127

$ws.Range("H45").Value = 1997.7273
$ws.Range("I45").Value = 1771.8889
$ws.Range("K45").Value = 1771.8889
$ws.Range("M45").Value = -1394.8889

$ws.Range("H75").Value = 80000
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0

$ws.Range("H76").Value = 64999.5
$ws.Range("J76").Value = 64999.5
$ws.Range("L76").Value = 64999.5
$ws.Range("N76").Value = -65675.5

$ws.Range("H78").Value = 80000
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0

$ws.Range("H79").Value = 64999.5
$ws.Range("J79").Value = 64999.5
$ws.Range("L79").Value = 64999.5
$ws.Range("N79").Value = -67339.5

$ws.Range("H97").Value = 1800.1111
$ws.Range("I97").Value = 1763.0454
$ws.Range("J97").Value = 1963.2
$ws.Range("K97").Value = 1763.0454
$ws.Range("L97").Value = 1963.2
$ws.Range("M97").Value = -1267.0454
$ws.Range("N97").Value = -2955.2

$ws.Range("H116").Value = 1728
$ws.Range("I116").Value = 1449.5
$ws.Range("K116").Value = 1449.5
$ws.Range("M116").Value = 844.5

$ws.Range("H122").Value = 2738.1794
$ws.Range("I122").Value = 1685.7826
$ws.Range("J122").Value = 4251
$ws.Range("K122").Value = 5057.3478
$ws.Range("L122").Value = 12753
$ws.Range("M122").Value = -2607.3478
$ws.Range("N122").Value = -17653

$ws.Range("H132").Value = 27781034
$ws.Range("I132").Value = 3355.9395
$ws.Range("J132").Value = 333335500
$ws.Range("K132").Value = 10067.8185
$ws.Range("L132").Value = 1000006500
$ws.Range("M132").Value = -7537.818499999999
$ws.Range("N132").Value = -1000011560

# --- Sheet BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1728
$ws.Range("I3").Value = 1449.5
$ws.Range("K3").Value = 1449.5
$ws.Range("M3").Value = -1335.5

$ws.Range("H94").Value = 911.25
$ws.Range("I94").Value = 722.37933
$ws.Range("J94").Value = 1693.7142
$ws.Range("K94").Value = 722.37933
$ws.Range("L94").Value = 1693.7142
$ws.Range("M94").Value = -271.37933
$ws.Range("N94").Value = -2595.7142

$ws.Range("H99").Value = 5076.25
$ws.Range("I99").Value = 3738
$ws.Range("J99").Value = 7752.75
$ws.Range("K99").Value = 3738
$ws.Range("L99").Value = 7752.75
$ws.Range("M99").Value = -2240
$ws.Range("N99").Value = -10748.75

$ws.Range("H105").Value = 18176.5
$ws.Range("I105").Value = 26214.75
$ws.Range("K105").Value = 26214.75
$ws.Range("M105").Value = -24467.75

$ws.Range("H107").Value = 3697.4443
$ws.Range("I107").Value = 3215.6
$ws.Range("J107").Value = 4299.75
$ws.Range("K107").Value = 3215.6
$ws.Range("L107").Value = 4299.75
$ws.Range("M107").Value = -1295.6
$ws.Range("N107").Value = -8139.75

# --- Sheet CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 25004160
$ws.Range("I31").Value = 3332.1924
$ws.Range("J31").Value = 71434270
$ws.Range("K31").Value = 3332.1924
$ws.Range("L31").Value = 71434270
$ws.Range("M31").Value = -3037.1924
$ws.Range("N31").Value = -71434860

$ws.Range("H34").Value = 25004160
$ws.Range("I34").Value = 3332.1924
$ws.Range("J34").Value = 71434270
$ws.Range("K34").Value = 3332.1924
$ws.Range("L34").Value = 71434270
$ws.Range("M34").Value = -3130.1924
$ws.Range("N34").Value = -71434674

$ws.Range("H58").Value = 2632.95
$ws.Range("I58").Value = 1977.3334
$ws.Range("J58").Value = 4599.8
$ws.Range("K58").Value = 1977.3334
$ws.Range("L58").Value = 4599.8
$ws.Range("M58").Value = -1774.3334
$ws.Range("N58").Value = -5005.8

$ws.Range("H107").Value = 1161.2
$ws.Range("I107").Value = 562.3
$ws.Range("K107").Value = 562.3
$ws.Range("M107").Value = 1357.7

$ws.Range("H134").Value = 1351.5555
$ws.Range("I134").Value = 1145.6875
$ws.Range("K134").Value = 3437.0625
$ws.Range("M134").Value = -902.0625

$ws.Range("H136").Value = 2632.95
$ws.Range("I136").Value = 1977.3334
$ws.Range("J136").Value = 4599.8
$ws.Range("K136").Value = 5932.0002
$ws.Range("L136").Value = 13799.4
$ws.Range("M136").Value = -3382.0002
$ws.Range("N136").Value = -18899.4

# --- Sheet CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H6").Value = 62
$ws.Range("I6").Value = 79.333336
$ws.Range("J6").Value = 10
$ws.Range("K6").Value = 238.000008
$ws.Range("L6").Value = 30
$ws.Range("M6").Value = -125.000008

$ws.Range("H40").Value = 184.1
$ws.Range("I40").Value = 77
$ws.Range("J40").Value = 344.75
$ws.Range("K40").Value = 308
$ws.Range("L40").Value = 1379
$ws.Range("M40").Value = -239
$ws.Range("N40").Value = -1517

$ws.Range("H130").Value = 2517.75
$ws.Range("I130").Value = 2015.5
$ws.Range("K130").Value = 6046.5

$ws.Range("H133").Value = 10432.695
$ws.Range("J133").Value = 15127.333
$ws.Range("L133").Value = 45381.999
$ws.Range("N133").Value = -55501.999

$ws.Range("H134").Value = 3811.5
$ws.Range("J134").Value = 14354.667
$ws.Range("L134").Value = 43064.001
$ws.Range("N134").Value = -53204.001

$ws.Range("H139").Value = 2645.7273
$ws.Range("J139").Value = 3271.25
$ws.Range("L139").Value = 9813.75
$ws.Range("N139").Value = -20093.75

$ws.Range("H140").Value = 1472.9445
$ws.Range("I140").Value = 887.3570999999999
$ws.Range("J140").Value = 3522.5
$ws.Range("K140").Value = 2662.0713
$ws.Range("L140").Value = 10567.5
$ws.Range("M140").Value = 2517.9287
$ws.Range("N140").Value = -20927.5

$ws.Range("H141").Value = 14301.733
$ws.Range("I141").Value = 16142.143
$ws.Range("J141").Value = 12691.375
$ws.Range("K141").Value = 48426.429
$ws.Range("L141").Value = 38074.125
$ws.Range("M141").Value = -43246.429
$ws.Range("N141").Value = -48434.125

# --- Sheet GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H64").Value = 44999
$ws.Range("I64").Value = 44999
$ws.Range("K64").Value = 44999
$ws.Range("M64").Value = -44751

$ws.Range("H67").Value = 44999
$ws.Range("I67").Value = 44999
$ws.Range("K67").Value = 44999
$ws.Range("M67").Value = -44141

$ws.Range("H80").Value = 3233.3333
$ws.Range("I80").Value = 3233.3333
$ws.Range("K80").Value = 3233.3333
$ws.Range("M80").Value = -2235.3333

$ws.Range("H83").Value = 3233.3333
$ws.Range("I83").Value = 3233.3333
$ws.Range("K83").Value = 16166.6665
$ws.Range("M83").Value = -11174.6665

$ws.Range("H97").Value = 1224
$ws.Range("I97").Value = 426.27274
$ws.Range("K97").Value = 426.27274
$ws.Range("M97").Value = 69.72726

$ws.Range("H103").Value = 43750
$ws.Range("J103").Value = 46666.668
$ws.Range("L103").Value = 46666.668
$ws.Range("N103").Value = -49010.668

$ws.Range("H122").Value = 7211.222
$ws.Range("I122").Value = 3320.3333
$ws.Range("K122").Value = 9960.999899999999
$ws.Range("M122").Value = -7510.999899999999

$ws.Range("H132").Value = 3227.7585
$ws.Range("I132").Value = 2933.762
$ws.Range("K132").Value = 8801.286
$ws.Range("M132").Value = -6271.286

# --- Sheet LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 3857.5
$ws.Range("J7").Value = 3805.5
$ws.Range("L7").Value = 3805.5
$ws.Range("N7").Value = -4029.5

$ws.Range("I16").Value = 689.3333
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 689.3333
$ws.Range("L16").Value = 200
$ws.Range("M16").Value = -519.3333
$ws.Range("N16").Value = -540

$ws.Range("H22").Value = 3027.611
$ws.Range("I22").Value = 2300
$ws.Range("K22").Value = 2300
$ws.Range("M22").Value = -2005

$ws.Range("H27").Value = 3027.611
$ws.Range("I27").Value = 2300
$ws.Range("K27").Value = 2300
$ws.Range("M27").Value = -2193

$ws.Range("H82").Value = 6441.7144
$ws.Range("I82").Value = 4049
$ws.Range("J82").Value = 7398.8
$ws.Range("K82").Value = 4049
$ws.Range("L82").Value = 7398.8
$ws.Range("M82").Value = -3688
$ws.Range("N82").Value = -8120.8

$ws.Range("H85").Value = 6441.7144
$ws.Range("I85").Value = 4049
$ws.Range("J85").Value = 7398.8
$ws.Range("K85").Value = 4049
$ws.Range("L85").Value = 7398.8
$ws.Range("M85").Value = -2801
$ws.Range("N85").Value = -9894.799999999999

$ws.Range("H93").Value = 2258.484
$ws.Range("I93").Value = 776.0714
$ws.Range("J93").Value = 3479.2942
$ws.Range("K93").Value = 776.0714
$ws.Range("L93").Value = 3479.2942
$ws.Range("M93").Value = 471.9286
$ws.Range("N93").Value = -5975.2942

$ws.Range("H115").Value = 75000
$ws.Range("J115").Value = 75000
$ws.Range("L115").Value = 75000

$ws.Range("H126").Value = 3857.5
$ws.Range("J126").Value = 3805.5
$ws.Range("L126").Value = 11416.5
$ws.Range("N126").Value = -16356.5

$ws.Range("H132").Value = 71431470
$ws.Range("I132").Value = 2786.0952
$ws.Range("K132").Value = 8358.285600000001
$ws.Range("M132").Value = -5828.285600000001

# --- Sheet WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H96").Value = 7543.1665
$ws.Range("I96").Value = 6474.5
$ws.Range("J96").Value = 8077.5
$ws.Range("K96").Value = 6474.5
$ws.Range("L96").Value = 8077.5
$ws.Range("M96").Value = -5101.5
$ws.Range("N96").Value = -10823.5

$ws.Range("H100").Value = 770776.9
$ws.Range("I100").Value = 1001565.1
$ws.Range("J100").Value = 1482.6666
$ws.Range("K100").Value = 2003130.2
$ws.Range("L100").Value = 2965.3332
$ws.Range("M100").Value = -2002589.2
$ws.Range("N100").Value = -4047.3332

$ws.Range("H113").Value = 594.4
$ws.Range("I113").Value = 232.10527
$ws.Range("K113").Value = 696.3158099999999
$ws.Range("M113").Value = 1473.68419

$ws.Range("H122").Value = 66668664
$ws.Range("J122").Value = 3277
$ws.Range("L122").Value = 9831
$ws.Range("N122").Value = -14731

$ws.Range("H124").Value = 42499.75
$ws.Range("J124").Value = 42499.75
$ws.Range("L124").Value = 42499.75
$ws.Range("N124").Value = -52319.75

$ws.Range("H126").Value = 3439.3125
$ws.Range("I126").Value = 3789.5715
$ws.Range("K126").Value = 11368.7145
$ws.Range("M126").Value = -8898.7145

$ws.Range("H132").Value = 3777
$ws.Range("I132").Value = 3922.2307
$ws.Range("J132").Value = 945
$ws.Range("K132").Value = 11766.6921
$ws.Range("L132").Value = 2835
$ws.Range("M132").Value = -9236.6921
$ws.Range("N132").Value = -7895

$ws.Range("H136").Value = 2652
$ws.Range("I136").Value = 2452.6365
$ws.Range("K136").Value = 7357.9095
$ws.Range("M136").Value = -4807.9095

# --- Additions of previously-absent cells ---
$wb.Worksheets.Item(5).Range("N6").Value = -256
$wb.Worksheets.Item(5).Range("M130").Value = -1026.5
$wb.Worksheets.Item(7).Range("N115").Value = -77350

# --- Deletions of cells (N75, N78 on ARM sheet) ---
$wb.Worksheets.Item(2).Range("N75").ClearContents()
$wb.Worksheets.Item(2).Range("N78").ClearContents()